# Add 3 denyeep column div areas to Image Freeway tab
# - Insert a new blank row at row 8 (pushes the existing (checkbox)/open-* rows
#   down by one, rows 8-16 -> 9-17)
# - Stamp a "///" divider row at row 6 (above the header row)
# - Stamp a "///" divider row at row 8 (the freshly inserted row, below the header)
# - Stamp a "///" divider row at row 18 (new last row, below the data)
# - Move the active selection to B14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row right after the header row (old row 7) so the
# "(checkbox)" / "open ..." rows shift from 8-16 down to 9-17.
$ws.Rows.Item(8).Insert()
# Inserting a row copies the formatting of the row above (the header's
# shaded style) onto the new row's cells; strip that back out so the new
# divider row starts from the default (unstyled) cell format.
$ws.Range("A8:I8").ClearFormats()

# Divider row above the header (row 6), columns A through I.
$topDivider = $ws.Range("A6:I6")
$topDivider.Value = "///"

# Divider row below the header (row 8, the newly inserted blank row),
# columns A through I.
$midDivider = $ws.Range("A8:I8")
$midDivider.Value = "///"

# Divider row at the very bottom (row 18), columns A through I.
$bottomDivider = $ws.Range("A18:I18")
$bottomDivider.Value = "///"

# Update the selection to match the saved workbook state.
$ws.Range("B14").Select()
